$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.478023052215576
$ws.Range("B1").Value = 1.971718192100525
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.601622700691223
$ws.Range("E1").Value = 0.6748588085174561
